$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: IBAN entry changes to "ES32 6784 345 0000" with commission 10
$ws.Range("A2").Value = "ES32 6784 345 0000"
$ws.Range("B2").Value = 10

# Update row 3: IBAN entry changes to "ES32 893 455 2333" with commission 7
$ws.Range("A3").Value = "ES32 893 455 2333"
$ws.Range("B3").Value = 7

# Delete rows 4 to 6 (old extra data no longer present)
$ws.Range("A4:C6").Delete() | Out-Null

# Move selection to the last data cell, matching the author's final cursor position
$ws.Range("C3").Select() | Out-Null
